$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 104.198
$ws.Range("C3").Value = 6.22
$ws.Range("C4").Value = 26.242
$ws.Range("C5").Value = 3.072
$ws.Range("C6").Value = 210.918
$ws.Range("C7").Value = 11.879
$ws.Range("C8").Value = 9.888999999999999
$ws.Range("C9").Value = 15.589
$ws.Range("C10").Value = 4.776
$ws.Range("C11").Value = 2.244
$ws.Range("C12").Value = 9.765000000000001
$ws.Range("C13").Value = 15.772
$ws.Range("C14").Value = 6.143
$ws.Range("C15").Value = 62.615
$ws.Range("V15").Value = "Green"
$ws.Range("C16").Value = 3.962
$ws.Range("V16").Value = "Green"
